# "I hope excel exportation is fully completed"
#
# Restructures the "Савдо" (Sales) sheet header row:
#   old: №, Категория, Махсулот, Сотилди, Жами сумма
#   new: №, Махсулот, Категория, Тан нархи, Нархи, Сотилди, Жами сумма
# (adds a "cost price" column and a "price" column, and swaps the
#  category/product column order), and moves the saved cell selections
# on the "Буюртмалар" and "Савдо" sheets.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Буюртмалар")
$ws2 = $wb.Worksheets.Item("Савдо")

# --- Sheet2 ("Савдо"): insert two new blank columns before the old D/E
#     (Сотилди / Жами сумма) columns, so they become F/G. ---
$ws2.Range("D1:E1").EntireColumn.Insert()

# Rebuild the header row with the new column order/content.
$ws2.Cells.Item(1, 2).Value = "Махсулот"     # B1 (was Категория)
$ws2.Cells.Item(1, 3).Value = "Категория"    # C1 (was Махсулот)
$ws2.Cells.Item(1, 4).Value = "Тан нархи"    # D1 (new)
$ws2.Cells.Item(1, 5).Value = "Нархи"        # E1 (new)
# F1/G1 already hold "Сотилди"/"Жами сумма" (shifted right by the insert).

# Match the new column widths to the "Категория"/"Махсулот" column width.
$ws2.Columns("D:E").ColumnWidth = $ws2.Columns("C").ColumnWidth()

# --- Saved selections ---
$ws2.Range("G8").Select()
$ws1.Range("E20").Select()
